$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.682.30'
$ws.Range("E2").Value = '  +5.65%  '
$ws.Range("D3").Value = '2.048.44'
$ws.Range("E3").Value = '  +3.13%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''251.78'
$ws.Range("E5").Value = '  +3.85%  '
$ws.Range("D6").Value = '''0.651'
$ws.Range("E6").Value = '  +1.68%  '
$ws.Range("D7").Value = '''65.14'
$ws.Range("E7").Value = '  +14.03%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").Value = '''59.95'
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").Value = '''0.378'
$ws.Range("E10").Value = '  +4.70%  '
$ws.Range("D11").Value = '''0.0761'
$ws.Range("E11").Value = '  +3.99%  '
$ws.Range("E12").Value = '  +1.68%  '
$ws.Range("D13").Value = '''0.920'
$ws.Range("E13").Value = '  -0.74%  '
$ws.Range("D14").Value = '''15.13'
$ws.Range("E14").Value = '  +7.08%  '
$ws.Range("D15").Value = '2.351.26'
$ws.Range("E15").Value = '  +3.34%  '
$ws.Range("D16").Value = '''20.72'
$ws.Range("E16").Value = '  +19.98%  '
$ws.Range("D17").Value = '''5.55'
$ws.Range("E17").Value = '  +5.45%  '
$ws.Range("D18").Value = '2.039.15'
$ws.Range("E18").Value = '  +2.68%  '
$ws.Range("D19").Value = '37.543.52'
$ws.Range("E19").Value = '  +5.62%  '
$ws.Range("D20").Value = '''74.01'
$ws.Range("E20").Value = '  +4.89%  '
$ws.Range("D21").Value = '0.0₃0875'
$ws.Range("E21").Value = '  +4.22%  '
$ws.Range("D22").Value = '''5.35'
$ws.Range("E22").Value = '  +5.39%  '
$ws.Range("D23").Value = '''238.63'
$ws.Range("E23").Value = '  +2.34%  '
$ws.Range("D24").Value = '''2.67'
$ws.Range("E24").Value = '  +13.96%  '
$ws.Range("E25").Value = '  -0.15%  '
$ws.Range("D26").Value = '''2.39'
$ws.Range("E26").Value = '  +3.82%  '
$ws.Range("D27").Value = '''9.59'
$ws.Range("E27").Value = '  +4.58%  '
$ws.Range("D28").Value = '''160.09'
$ws.Range("E28").Value = '  -2.11%  '
$ws.Range("D29").Value = '''19.94'
$ws.Range("E29").Value = '  +2.14%  '
$ws.Range("E30").Value = '  +2.35%  '
$ws.Range("E31").Value = '  +26.51%  '
$ws.Range("D32").Value = '''5.21'
$ws.Range("E32").Value = '  +8.30%  '
$ws.Range("D33").Value = '''1.20'
$ws.Range("E33").Value = '  +6.27%  '
$ws.Range("D34").Value = '''4.72'
$ws.Range("E34").Value = '  +10.60%  '
$ws.Range("D35").Value = '''0.0618'
$ws.Range("E35").Value = '  +4.89%  '
$ws.Range("E36").Value = '  +1.15%  '
$ws.Range("D37").Value = '''1.86'
$ws.Range("E37").Value = '  +3.04%  '
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("D39").Value = '''6.08'
$ws.Range("E39").Value = '  +23.12%  '
$ws.Range("D40").Value = '''0.104'
$ws.Range("E40").Value = '  +16.32%  '
$ws.Range("D41").Value = '''2.82'
$ws.Range("E41").Value = '  +24.69%  '
$ws.Range("D42").Value = '''1.23'
$ws.Range("E42").Value = '  +3.46%  '
$ws.Range("E43").Value = '  +4.28%  '
$ws.Range("E44").Value = '  +3.69%  '
$ws.Range("E45").Value = '  +5.09%  '
$ws.Range("D46").Value = '''16.96'
$ws.Range("E46").Value = '  +9.60%  '
$ws.Range("D47").Value = '''7.99'
$ws.Range("E47").Value = '  +7.51%  '
$ws.Range("D48").Value = '''95.05'
$ws.Range("E48").Value = '  +4.38%  '
$ws.Range("D49").Value = '1.413.66'
$ws.Range("E49").Value = '  +2.22%  '
$ws.Range("E50").Value = '  +2.21%  '
$ws.Range("D51").Value = '''47.24'
$ws.Range("E51").Value = '  +2.95%  '
